$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6 (Leve Item ID 4564): 
$ws.Cells.Item(6, 8).Value = 378.1875
$ws.Cells.Item(6, 9).Value = 86.454544
$ws.Cells.Item(6, 10).Value = 1020
$ws.Cells.Item(6, 11).Value = 259.363632
$ws.Cells.Item(6, 12).Value = 3060
$ws.Cells.Item(6, 13).Value = -147.363632
$ws.Cells.Item(6, 14).Value = -3284

# Row 31 (Leve Item ID 4576): 
$ws.Cells.Item(31, 8).Value = 1701.8
$ws.Cells.Item(31, 9).Value = 127.25
$ws.Cells.Item(31, 11).Value = 381.75
$ws.Cells.Item(31, 13).Value = -151.75

# Row 88 (Leve Item ID 12608): 
$ws.Cells.Item(88, 8).Value = 4736.64
$ws.Cells.Item(88, 9).Value = 2406.889
$ws.Cells.Item(88, 10).Value = 6047.125
$ws.Cells.Item(88, 11).Value = 2406.889
$ws.Cells.Item(88, 12).Value = 6047.125
$ws.Cells.Item(88, 13).Value = -2000.889
$ws.Cells.Item(88, 14).Value = -6859.125

# Row 91 (Leve Item ID 12608): 
$ws.Cells.Item(91, 8).Value = 4736.64
$ws.Cells.Item(91, 9).Value = 2406.889
$ws.Cells.Item(91, 10).Value = 6047.125
$ws.Cells.Item(91, 11).Value = 2406.889
$ws.Cells.Item(91, 12).Value = 6047.125
$ws.Cells.Item(91, 13).Value = -1002.889
$ws.Cells.Item(91, 14).Value = -8855.125

$ws = $wb.Worksheets.Item("ARM")
# Row 6 (Leve Item ID 2226): 
$ws.Cells.Item(6, 8).Value = 31150.4
$ws.Cells.Item(6, 9).Value = 37688.5
$ws.Cells.Item(6, 11).Value = 37688.5
$ws.Cells.Item(6, 13).Value = -37515.5

# Row 97 (Leve Item ID 19941): 
$ws.Cells.Item(97, 8).Value = 776.5714
$ws.Cells.Item(97, 9).Value = 770
$ws.Cells.Item(97, 10).Value = 785.3333
$ws.Cells.Item(97, 11).Value = 770
$ws.Cells.Item(97, 12).Value = 785.3333
$ws.Cells.Item(97, 13).Value = -274
$ws.Cells.Item(97, 14).Value = -1777.3333

# Row 102 (Leve Item ID 19945): 
$ws.Cells.Item(102, 8).Value = 3616.6667
$ws.Cells.Item(102, 10).Value = 3925
$ws.Cells.Item(102, 12).Value = 3925
$ws.Cells.Item(102, 14).Value = -7169

$ws = $wb.Worksheets.Item("BSM")
# Row 20 (Leve Item ID 14149): 
$ws.Cells.Item(20, 8).Value = 1457.5927
$ws.Cells.Item(20, 9).Value = 1329.0667
$ws.Cells.Item(20, 10).Value = 1618.25
$ws.Cells.Item(20, 11).Value = 1329.0667
$ws.Cells.Item(20, 12).Value = 1618.25
$ws.Cells.Item(20, 13).Value = -1082.0667
$ws.Cells.Item(20, 14).Value = -2112.25

# Row 62 (Leve Item ID 10586): 
$ws.Cells.Item(62, 8).Value = 30000
$ws.Cells.Item(62, 10).Value = 30000
$ws.Cells.Item(62, 12).Value = 30000
$ws.Cells.Item(62, 14).Value = -31372

# Row 65 (Leve Item ID 10586): 
$ws.Cells.Item(65, 8).Value = 30000
$ws.Cells.Item(65, 10).Value = 30000
$ws.Cells.Item(65, 12).Value = 90000
$ws.Cells.Item(65, 14).Value = -96864

# Row 86 (Leve Item ID 12526): 
$ws.Cells.Item(86, 8).Value = 1493.6046
$ws.Cells.Item(86, 9).Value = 1414.7368
$ws.Cells.Item(86, 10).Value = 1556.0416
$ws.Cells.Item(86, 11).Value = 1414.7368
$ws.Cells.Item(86, 12).Value = 1556.0416
$ws.Cells.Item(86, 13).Value = -291.7367999999999
$ws.Cells.Item(86, 14).Value = -3802.0416

# Row 89 (Leve Item ID 12526): 
$ws.Cells.Item(89, 8).Value = 1493.6046
$ws.Cells.Item(89, 9).Value = 1414.7368
$ws.Cells.Item(89, 10).Value = 1556.0416
$ws.Cells.Item(89, 11).Value = 7073.683999999999
$ws.Cells.Item(89, 12).Value = 7780.208000000001
$ws.Cells.Item(89, 13).Value = -1457.683999999999
$ws.Cells.Item(89, 14).Value = -19012.208

# Row 94 (Leve Item ID 19939): 
$ws.Cells.Item(94, 8).Value = 991.1429000000001
$ws.Cells.Item(94, 9).Value = 779.63635
$ws.Cells.Item(94, 11).Value = 779.63635
$ws.Cells.Item(94, 13).Value = -328.63635

# Row 99 (Leve Item ID 19943): 
$ws.Cells.Item(99, 8).Value = 2402.7222
$ws.Cells.Item(99, 9).Value = 1987.5
$ws.Cells.Item(99, 10).Value = 2521.3572
$ws.Cells.Item(99, 11).Value = 1987.5
$ws.Cells.Item(99, 12).Value = 2521.3572
$ws.Cells.Item(99, 13).Value = -489.5
$ws.Cells.Item(99, 14).Value = -5517.3572

# Row 105 (Leve Item ID 19947): 
$ws.Cells.Item(105, 8).Value = 2199.35
$ws.Cells.Item(105, 9).Value = 2177.2407
$ws.Cells.Item(105, 10).Value = 2398.3333
$ws.Cells.Item(105, 11).Value = 2177.2407
$ws.Cells.Item(105, 12).Value = 2398.3333
$ws.Cells.Item(105, 13).Value = -430.2406999999998
$ws.Cells.Item(105, 14).Value = -5892.3333

$ws = $wb.Worksheets.Item("CRP")
# Row 134 (Leve Item ID 44020): 
$ws.Cells.Item(134, 8).Value = 3307.353
$ws.Cells.Item(134, 9).Value = 971.2
$ws.Cells.Item(134, 10).Value = 6644.7144
$ws.Cells.Item(134, 11).Value = 2913.6
$ws.Cells.Item(134, 12).Value = 19934.1432
$ws.Cells.Item(134, 13).Value = -378.6000000000004
$ws.Cells.Item(134, 14).Value = -25004.1432

# Row 141 (Leve Item ID 43345): 
$ws.Cells.Item(141, 8).Value = 47142.625
$ws.Cells.Item(141, 10).Value = 47142.625
$ws.Cells.Item(141, 12).Value = 47142.625
$ws.Cells.Item(141, 14).Value = -57502.625

$ws = $wb.Worksheets.Item("CUL")
# Row 131 (Leve Item ID 36060): 
$ws.Cells.Item(131, 8).Value = 922.41
$ws.Cells.Item(131, 9).Value = 499.66666
$ws.Cells.Item(131, 10).Value = 935.48456
$ws.Cells.Item(131, 11).Value = 1498.99998
$ws.Cells.Item(131, 12).Value = 2806.45368
$ws.Cells.Item(131, 13).Value = 3541.00002
$ws.Cells.Item(131, 14).Value = -12886.45368

$ws = $wb.Worksheets.Item("GSM")
# Row 80 (Leve Item ID 12521): 
$ws.Cells.Item(80, 8).Value = 2862.5
$ws.Cells.Item(80, 9).Value = 2775
$ws.Cells.Item(80, 10).Value = 2950
$ws.Cells.Item(80, 11).Value = 2775
$ws.Cells.Item(80, 12).Value = 2950
$ws.Cells.Item(80, 13).Value = -1777
$ws.Cells.Item(80, 14).Value = -4946

# Row 83 (Leve Item ID 12521): 
$ws.Cells.Item(83, 8).Value = 2862.5
$ws.Cells.Item(83, 9).Value = 2775
$ws.Cells.Item(83, 10).Value = 2950
$ws.Cells.Item(83, 11).Value = 13875
$ws.Cells.Item(83, 12).Value = 14750
$ws.Cells.Item(83, 13).Value = -8883
$ws.Cells.Item(83, 14).Value = -24734

# Row 97 (Leve Item ID 19940): 
$ws.Cells.Item(97, 8).Value = 910

$ws = $wb.Worksheets.Item("LTW")
# Row 55 (Leve Item ID 5284): 
$ws.Cells.Item(55, 8).Value = 614.2727
$ws.Cells.Item(55, 9).Value = 274.42856
$ws.Cells.Item(55, 11).Value = 274.42856
$ws.Cells.Item(55, 13).Value = -101.42856

# Row 68 (Leve Item ID 12563): 
$ws.Cells.Item(68, 8).Value = 3100
$ws.Cells.Item(68, 9).Value = 0
$ws.Cells.Item(68, 11).Value = 0
$ws.Cells.Item(68, 13).ClearContents()

# Row 71 (Leve Item ID 12563): 
$ws.Cells.Item(71, 8).Value = 3100
$ws.Cells.Item(71, 9).Value = 0
$ws.Cells.Item(71, 11).Value = 0
$ws.Cells.Item(71, 13).ClearContents()

# Row 82 (Leve Item ID 12565): 
$ws.Cells.Item(82, 8).Value = 1085.3334
$ws.Cells.Item(82, 9).Value = 1000
$ws.Cells.Item(82, 10).Value = 1091.4286
$ws.Cells.Item(82, 11).Value = 1000
$ws.Cells.Item(82, 12).Value = 1091.4286
$ws.Cells.Item(82, 13).Value = -639
$ws.Cells.Item(82, 14).Value = -1813.4286

# Row 85 (Leve Item ID 12565): 
$ws.Cells.Item(85, 8).Value = 1085.3334
$ws.Cells.Item(85, 9).Value = 1000
$ws.Cells.Item(85, 10).Value = 1091.4286
$ws.Cells.Item(85, 11).Value = 1000
$ws.Cells.Item(85, 12).Value = 1091.4286
$ws.Cells.Item(85, 13).Value = 248
$ws.Cells.Item(85, 14).Value = -3587.4286

# Row 93 (Leve Item ID 19993): 
$ws.Cells.Item(93, 8).Value = 29254.8
$ws.Cells.Item(93, 9).Value = 854.1429000000001
$ws.Cells.Item(93, 10).Value = 65401.09
$ws.Cells.Item(93, 11).Value = 854.1429000000001
$ws.Cells.Item(93, 12).Value = 65401.09
$ws.Cells.Item(93, 13).Value = 393.8570999999999
$ws.Cells.Item(93, 14).Value = -67897.09

# Row 100 (Leve Item ID 19995): 
$ws.Cells.Item(100, 8).Value = 1525.4166
$ws.Cells.Item(100, 9).Value = 1650
$ws.Cells.Item(100, 10).Value = 1500.5
$ws.Cells.Item(100, 11).Value = 1650
$ws.Cells.Item(100, 12).Value = 1500.5
$ws.Cells.Item(100, 13).Value = -1109
$ws.Cells.Item(100, 14).Value = -2582.5

$ws = $wb.Worksheets.Item("WVR")
# Row 81 (Leve Item ID 12596): 
$ws.Cells.Item(81, 8).Value = 1180
$ws.Cells.Item(81, 9).Value = 1180
$ws.Cells.Item(81, 11).Value = 2360
$ws.Cells.Item(81, 13).Value = -1299

# Row 84 (Leve Item ID 12596): 
$ws.Cells.Item(84, 8).Value = 1180
$ws.Cells.Item(84, 9).Value = 1180
$ws.Cells.Item(84, 11).Value = 11800
$ws.Cells.Item(84, 13).Value = -6496

# Row 96 (Leve Item ID 19977): 
$ws.Cells.Item(96, 8).Value = 4334.6665
$ws.Cells.Item(96, 9).Value = 2000
$ws.Cells.Item(96, 10).Value = 5502
$ws.Cells.Item(96, 11).Value = 2000
$ws.Cells.Item(96, 12).Value = 5502
$ws.Cells.Item(96, 13).Value = -627
$ws.Cells.Item(96, 14).Value = -8248
